$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.126.20'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = '3.413.48'
$ws.Range("E3").Value = '  -1.27%  '
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.35'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.08'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("E7").Value = '  +3.79%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '3.416.02'
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.131'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.70%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.97'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.411'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.93%  '
$ws.Range("D13").Value = '4.009.41'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("E14").Value = '  +0.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.07'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.11%  '
$ws.Range("D16").Value = '66.188.50'
$ws.Range("E16").Value = '  -0.59%  '
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").Value = '3.415.55'
$ws.Range("E18").Value = '  -1.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.81%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '365.41'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.53'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.83%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.998'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.533'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("B26").Value = 'PEPE'
$ws.Range("C26").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000125'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.53%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.78'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.41%  '
$ws.Range("E28").Value = '  +1.33%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.98'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.57%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.73'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '23.28'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.20%  '
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.25'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.96'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.44%  '
$ws.Range("E38").Value = '  -2.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.78'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.56'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.32%  '
$ws.Range("D43").Value = '2.693.87'
$ws.Range("E43").Value = '  -1.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0686'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.96%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.08'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '39.91'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.98%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '334.40'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.80%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0285'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.96%  '
$ws.Range("E50").Value = '  +2.50%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '31.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.83%  '
